$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RM 8): F was missing -> now has a value
$ws.Range("F3").Value = 17.64

# Row 5 (RM 14): F had a value -> now missing
$ws.Range("F5").ClearContents()

# Row 21 (RM 135): F was missing -> now has a value
$ws.Range("F21").Value = 16.58

# Row 23 (RM 140): F had a value -> now missing
$ws.Range("F23").ClearContents()

# Remove the "RM 232" row (26) and the "SC 92" row (now at 27 after the
# first delete) so the remaining SC rows shift up two positions.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Row 32 (SC 193) now has a value for F where it used to be missing.
$ws.Range("F32").Value = 17.39
